# Update automatico via Actualizar 02-12-2021 16-21-59
#
# Column D holds a "last updated" timestamp for each of the 14 logical
# rows (r2..r15). Each refresh shifts the previous timestamps down into
# the next block of 14 rows, dropping the oldest block, and stamps the
# newest timestamp into the first block.
#
#   rows  2..15 (newest) <- new timestamp 44239.68184616481
#   rows 16..29           <- old value that was in rows  2..15 (44239.66063602377 ~= 44239.66063601852)
#   rows 30..43 (oldest)  <- old value that was in rows 16..29 (44239.63938445602)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp   = 44239.68184616481
$shiftedBlock1  = 44239.66063601852
$shiftedBlock2  = 44239.63938445602

# rows 30..43 (oldest block) take the timestamp previously held by rows 16..29
for ($row = 30; $row -le 43; $row++) {
    $ws.Cells.Item($row, 4).Value2 = $shiftedBlock2
}

# rows 16..29 take the timestamp previously held by rows 2..15
for ($row = 16; $row -le 29; $row++) {
    $ws.Cells.Item($row, 4).Value2 = $shiftedBlock1
}

# rows 2..15 (newest block) get stamped with the brand new update timestamp
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 4).Value2 = $newTimestamp
}
